$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Tighten the decimal number format used by column B (numFmtId 165)
#    from "0.000" to "###0.000"
$ws.Range("B1:B10").NumberFormat = "###0.000"

# 3. Update the cell A1 comment text (CEIC Data Manager metadata blob)
$comment = $ws.Range("A1").Comment
$comment.Text("GRoAAB+LCAAAAAAAAAOlGdtu48b1Vwg9tUAlkpK9sY1ZBrrZEaqLIcn1Oi/FiBxLrKmhQg4t6y0FEqRIUxRFsSnSK/qUokC3izYB0t1e/iVYe7dP/YWeuXBISnK22i6MNefc5syZcx2jt2/mgXFNotgP6cOSXbFKBqFu6Pl0+rCUsMuy/aD0toPaNy4JTnGE54QBsQFcND66if2HpRljiyPTXC6XlWWtEkZTs2pZtvmo1x25MzLHZZ/GDFOXlDSX93qukoOa3rxHGPYww5LzYakz6lSaxHdbAOthiqckqjSS2KckjtuU+cwnMeeMCGak2ep9Tx7MqVYeVGxkbsAzykbiB56kK1BKuKKDbcnYnxOnatkHZeugXLPH1uGRvX9Uq1bsg713U0ZNiLo4ZiMSXfuuAIwYni8Eu3VQs61De79WReZWIpCVGcBBg8Abkms/Jl6TBEG8k0VMdYF1l8GpdzOmhcwcrxL05iqcRHgxG/ssILupMew1jDlVumRCHHQcRsQF+72RSn2yHETKrONFF7DjmR+xVQuvdpZ1FpNosOBG2o3VQa2QsnpAIna2gLsmHrgCIBwWJQSZ9yAzppYfu/Dt04R4ziUO4jxTAYnOw+gqXmCX9CGOTS5jSYMQe+BwzI+Z78aZgA0MOo3CBYiE3Rth4B2D2JR6C0bL7lAwMt+4EYZX6+KLSCTuVdww3Ooca/kbcDSahcsBDVajZBK7kT8hXquRUm/FIR6SiruZxCycgxYZCElYDrKCfxCC62DUIq4/x8FpAGaMnRpIKQBQPWHhpc+aYZDMqbbnGhSdw4nG5EafUK/RAK6XcquHtENTeukKW1FFhmG41HtuIoQRcuB67Oob2cSsU7cApnxyEyFuhJ/y2A+gQuTvIgctesVoRgjb6hISg3gyPOY1x2ms+sl8AgE2gSi7FrvGyMzwCBwVvB30ciyoI2XxM7asI/EDemg0alPvfroUiWC73F6ODbg1EIIzBY0A0yuAnvts1q+nZ9mCQdIC99Jv4hDE7iLAKwHWVsrDUIe6QeIRmRI69FK4KNdNUd+LRxugLkS5gzBdjVcLnhrMeyha5BInAVQrBkllmuWLNTCqx1frNHkQOouC9MYd3gvE0Ay43rziQrrkBa/ihnMOMKEGn4+QmafnOd8lbTrtYjpNIKtqO67Dtb/xfDCOMI35cXQSXXO97UQovReZ3B15WYNEuJO8rBCwyFyjQ2MyX4QRDnpgGP84oaKAqgoBFu5hNlMriOWAuKmRzYxVcxU1SxV/HZkISnkMnvHzQZmDCip+GNl3yBAvwhA/Zi/0oPLjwJ9EBS/bioMby+ph6m/8dDvWxvQSoPeEdPNdsuINSbZQcO60jp0ixAo1QSVnNNw7qO5bvMkSayROPCQ4MNrQlTJidOg1idkc2I6MIYl9D758HBwZ75AJ8SHqhYlU3O3MnedDxxF5L4HueiVUqYO+RUiRABLs1Kc42CTUmIzBuSA4ClY5QnnUbugC3d2P/nX766cvnn1+9/HjV19+8J+///LFP356++RD+Lj7y19vP/mFPKYkRmM8CYhQaNw4OLBqe+BoGoRUfoCS7yUuE7CLC9EC6DVSvaxYNNud5km3IRKKBqbsvL1IeDB08SpMsuVIHkJsJK7UTD1BkjjjNEGpdQHbIrzuyyjkbes1KVLn8fcxSlu8fP75y+d/updbGSwrLvbh4X7Zrr629sAAYG/Q6drDxwHZ7aXDxl7Z2i9XqzniNRo0hJkJmjttp47n8PnCqtYsWydzTzvyNqJ1lJI0xlNzjU+CmmFCWbTSLpBfp0jh+GMIEY2WoZBbKBf98sev/vy4QKWsqyBFKaBcmETSm8x0IUT3h2NjNDgbNtvGuD3ifpLhcnRS+DcQq911PBWcitIEB9+BSZnPyUYJWsCSEV4aBLszYwWRmIvDgrNtg8qN3lDkupYnUZgs5I3kGDLoFkqdTbZybMk1AifsuZF0MtQWcqnr7d++2MagDqLc7Iz6TI98eRgqYCQoh1dR+9k/X3z10Ytnz+6e/uz2qx8WJKh99CQEfg7RlF9qt4eUp+rNGgSdj4Qxr6zv5+qLAvLO8TT0KYsd+4FoGtUKAavNpYnfqDOHkicEC3sBfA2C3sFx+4apwHb6yCwCQM8FhmobZs22Bsgcntn137/57d2vvrj79Omrj/54+/Efbj/59OXz37168nsZdXePn9795InK8uuFQOjCW3jZBRpiInQNHo0Gr93G1+//3KAhM6DnMBKRkb5+/7OcMK6o6E4yydDTaUWKKmyQ5pk5n5FTRetQ4NMssgFo8hJW0xSqiIUL3802ebfMRfG4E4hvdcblJCZGCO3Ut+EkReKM+X/lUyyypJ6+ZVXtqsJKbfgRJjjOmf4kCCfQZKQIMXGtkRS4vpkhoxX7nXQHjXo3I5FKDCIPBjSLj438A6U9JS8pnThd6ckigwAWGj83CfiQvEG2idKSc2nMVMPmZd3j6W/7bFaggIE9imRDRNXz5ShZQDfMVAN7P1681OT6377sVfMdcbbutIp4WOewUAiLaA4QeJGaFEqmqU7MB1vZzva5abIl4AqvO2AO9UIpO61r6Csjk+eddhSF0dbkk2FSsh500pBRzMzimkbcqey6veyuUkCa8PSHHP3UCcMWCQjb7fnOzLh74fUb88Ld78raiQeBp4y52+ihzZIJyL9hckf5f58wpbPVowgaK/7isfObYzq5DmHg3VEbeRTByCdA2F09Ch77Ucwe8UygviTkQkMuZIf6iA9c8kOsL5xDuQa8mRduFrRMI5fJh+Qw6Ppzf8ep0ErDuygETLlYyA6us5uj8MrSJzfQX+YkQE6c/ACqBp94dpMm/RVSqebnbzWxP52xXRV7a4KJRyZW2Z2QannPsw7Kh4TUyrYN/2O3WrWsff7Qo4RD4vDJcsdNzPTCsj/wOP8FUaK7XRkaAAA=")
